# Add "type" column (M) to the accounts table, populate a second account
# row (row 3) as a copy of the first account row (row 2), and tag it as
# "varified" in the new type column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for column M
$ws.Range("M1").Value = "type"

# Duplicate the first account (row 2) into row 3, values + formatting
$ws.Range("A2:H2").Copy($ws.Range("A3:H3"))

# Mark the new account row as "varified" in the type column
$ws.Range("M3").Value = "varified"

# Match the author's final selection
$ws.Range("N7").Select()
